$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Clear out the old B2:B5 question cells (their text is being relocated / reworded
# further down the sheet). B1 (rich-text "pixel spacing..." header) is left untouched.
$ws.Range("B2").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()

# Previously-empty rows further down now hold the (moved/reworded) questions.
# NOTE: order matters for how new shared strings get appended/indexed, so we set
# these in the same first-seen order as the final sharedStrings table.
$ws.Range("B10").Value = "要不要先做分类分析 然后再具体分析+框"
$ws.Range("B13").Value = "activation function等之类的需要什么注意的? (sigmoid, Relu,tanh等?)"

# B3 keeps a value, but the wording changed and it also got expanded with more detail.
$ws.Range("B3").Value = "图像需要做什么样的data preprocessing吗? 例如color standardization, b/w pixel rate, pic rotation, pic zoom in/out; pixel normalization之类的?"

$ws.Range("B14").Value = "CNN推荐的python的包? Keras,mxnet,sklearn-theano,Mask-RCNN等"

# New discussion notes appended below the numbered list (rows 21-26).
$ws.Range("B21").Value = "敏感区域系数百分比"
$ws.Range("B22").Value = "整个图片的比较"
$ws.Range("B23").Value = "keras比较好上手 但是比较简练"
$ws.Range("B24").Value = "tensorflow"
$ws.Range("B25").Value = "整个图的分割 vs 权重的输出 vs 其他的输出"
$ws.Range("B26").Value = "luna的肺结节比赛?"

# This one reuses an already-existing shared string, so its position doesn't affect
# new-string indexing; set it last to mirror the move described in the diff.
$ws.Range("B12").Value = "一般图像分析的方法:CNN,DNN?"

# Update the page setup (matches paperSize/orientation seen in the saved file).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Scroll/selection state: window now shows row 2 at the top with D19 selected.
$win = $excel.ActiveWindow
$ws.Range("D19").Select()
$win.ScrollRow = 2
$win.ScrollColumn = 1
